# GDE-9324 - addressed comments
# Rename the "Clients" masterlist sheet/column to the new UAT Deal Scenario
# naming convention and refresh the related view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet rename: "Clients" -> "UAT_Deal_Scenarios"
$ws.Name = "UAT_Deal_Scenarios"

# 2) Header rename: "UAT_Client" -> "UAT_Deal_Scenario_Name"
$ws.Range("B1").Value = "UAT_Deal_Scenario_Name"

# 3) Column B widens (bestFit) to accommodate the longer header text
$ws.Columns.Item(2).ColumnWidth = 25

# 4) Active cell/selection moves to B9
$ws.Range("B9").Select()

# 5) Restore the workbook window position recorded at save time
$excel.Windows.Item(1).Left = -1875
$excel.Windows.Item(1).Top = 5415
